# se modifica data para hacer prueba en QA
$wb = $excel.ActiveWorkbook

# --- DatosCuenta sheet ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokQATres"
$wsCuenta.Range("B2").Value = "SmokeNameQATres"
$wsCuenta.Range("C2").Value = 27100119
$wsCuenta.Range("D2").Value = 121

# --- DatosHogar sheet ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 640

# --- DatosAP sheet ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200122

$wb.Save()
